$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1:E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
